# cambios de may de mayo
$wb = $excel.ActiveWorkbook

# --- Worksheet "Reporte de Formatos": update row 8 values (shift report period
#     from Q4 2021 to Q1 2022) ---
$ws = $wb.Worksheets.Item("Reporte de Formatos")

$ws.Range("A8").Value = 2022
$ws.Range("B8").Value = 44562
$ws.Range("C8").Value = 44651
$ws.Range("N8").Value = 44659
$ws.Range("O8").Value = 44659

# --- Sheet view / selection changes (scroll so column N is at the left edge,
#     select R10) ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$win.ScrollRow = 2
$ws.Range("R10").Select()

$wb.Save()
